$wb = $excel.ActiveWorkbook
$maxes = $wb.Worksheets.Item("Maxes")
$upper1 = $wb.Worksheets.Item("Upper1")
$theory = $wb.Worksheets.Item("Theoretical Weight Scheme")

$maxes.Range("F1:G12").Copy($theory.Range("A1"))
$theory.Range("A1:B12").RowHeight = 15.75

$upper1.Range("F12").Select()
$theory.Range("D17").Select()

$maxes.Columns.Item(6).ColumnWidth = 16.5
$maxes.Select()
$maxes.Range("F22").Select()
